$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2292.2
$ws.Range("I38").Value = 938.4
$ws.Range("J38").Value = 4999.8
$ws.Range("K38").Value = 2815.2
$ws.Range("L38").Value = 14999.4
$ws.Range("M38").Value = -2443.2
$ws.Range("N38").Value = -15743.4

$ws.Range("H52").Value = 1388.1
$ws.Range("I52").Value = 431.22223
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 1293.66669
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -1133.66669
$ws.Range("N52").Value = -30320

$ws.Range("H69").Value = 8598
$ws.Range("I69").Value = 6003.6665
$ws.Range("J69").Value = 12489.5
$ws.Range("K69").Value = 18010.9995
$ws.Range("L69").Value = 37468.5
$ws.Range("M69").Value = -17136.9995
$ws.Range("N69").Value = -39216.5

$ws.Range("H72").Value = 8598
$ws.Range("I72").Value = 6003.6665
$ws.Range("J72").Value = 12489.5
$ws.Range("K72").Value = 54032.9985
$ws.Range("L72").Value = 112405.5
$ws.Range("M72").Value = -49664.9985
$ws.Range("N72").Value = -121141.5

$ws.Range("H93").Value = 400000
$ws.Range("J93").Value = 400000
$ws.Range("L93").Value = 400000
$ws.Range("N93").Value = -404992

$ws.Range("H112").Value = 973.5
$ws.Range("J112").Value = 1161.25
$ws.Range("L112").Value = 3483.75
$ws.Range("N112").Value = -5699.75

$ws.Range("H121").Value = 1685
$ws.Range("J121").Value = 1685
$ws.Range("L121").Value = 5055
$ws.Range("N121").Value = -8549

$ws.Range("H137").Value = 1764.3103
$ws.Range("I137").Value = 1763.24
$ws.Range("K137").Value = 5289.72
$ws.Range("M137").Value = -2739.72

$ws.Range("H138").Value = 1748.4822
$ws.Range("J138").Value = 2166.7632
$ws.Range("L138").Value = 6500.2896
$ws.Range("N138").Value = -16780.2896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 19510.223
$ws.Range("I43").Value = 10342
$ws.Range("J43").Value = 20656.25
$ws.Range("K43").Value = 10342
$ws.Range("L43").Value = 20656.25
$ws.Range("M43").Value = -10029
$ws.Range("N43").Value = -21282.25

$ws.Range("H46").Value = 4569
$ws.Range("I46").Value = 3138
$ws.Range("K46").Value = 3138
$ws.Range("M46").Value = -2819

$ws.Range("H61").Value = 4605.0435
$ws.Range("J61").Value = 4486
$ws.Range("L61").Value = 4486
$ws.Range("N61").Value = -4910

$ws.Range("H74").Value = 2344.8718
$ws.Range("I74").Value = 1757.4138
$ws.Range("J74").Value = 4048.5
$ws.Range("K74").Value = 1757.4138
$ws.Range("L74").Value = 4048.5
$ws.Range("M74").Value = -883.4138
$ws.Range("N74").Value = -5796.5

$ws.Range("H77").Value = 2344.8718
$ws.Range("I77").Value = 1757.4138
$ws.Range("J77").Value = 4048.5
$ws.Range("K77").Value = 8787.069
$ws.Range("L77").Value = 20242.5
$ws.Range("M77").Value = -4419.069
$ws.Range("N77").Value = -28978.5

$ws.Range("H110").Value = 1378
$ws.Range("J110").Value = 2337.6
$ws.Range("L110").Value = 2337.6
$ws.Range("N110").Value = -6427.6

$ws.Range("H122").Value = 3472.8333
$ws.Range("I122").Value = 3324.8
$ws.Range("J122").Value = 3842.9167
$ws.Range("K122").Value = 9974.400000000001
$ws.Range("L122").Value = 11528.7501
$ws.Range("M122").Value = -7524.400000000001
$ws.Range("N122").Value = -16428.7501

$ws.Range("H136").Value = 4605.0435
$ws.Range("J136").Value = 4486
$ws.Range("L136").Value = 13458
$ws.Range("N136").Value = -18558

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5050
$ws.Range("I41").Value = 5050
$ws.Range("K41").Value = 5050
$ws.Range("M41").Value = -4622

$ws.Range("H58").Value = 2660.3076
$ws.Range("I58").Value = 2681.889
$ws.Range("K58").Value = 2681.889
$ws.Range("M58").Value = -2478.889

$ws.Range("H99").Value = 3226.7576
$ws.Range("J99").Value = 3456.5454
$ws.Range("L99").Value = 3456.5454
$ws.Range("N99").Value = -6452.5454

$ws.Range("H107").Value = 1467.4
$ws.Range("I107").Value = 1011
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1011
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 909
$ws.Range("N107").Value = -5340

$ws.Range("H126").Value = 3226.7576
$ws.Range("J126").Value = 3456.5454
$ws.Range("L126").Value = 10369.6362
$ws.Range("N126").Value = -15309.6362

$ws.Range("H134").Value = 1868.9333
$ws.Range("I134").Value = 1784.5652
$ws.Range("K134").Value = 5353.6956
$ws.Range("M134").Value = -2818.6956

$ws.Range("H136").Value = 2660.3076
$ws.Range("I136").Value = 2681.889
$ws.Range("K136").Value = 8045.667
$ws.Range("M136").Value = -5495.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 3333.6667
$ws.Range("J41").Value = 3501
$ws.Range("L41").Value = 10503
$ws.Range("N41").Value = -11179

$ws.Range("H54").Value = 5666
$ws.Range("J54").Value = 6499
$ws.Range("L54").Value = 19497
$ws.Range("N54").Value = -20615

$ws.Range("H59").Value = 6827.1816
$ws.Range("I59").Value = 3749.5
$ws.Range("J59").Value = 7511.1113
$ws.Range("K59").Value = 11248.5
$ws.Range("L59").Value = 22533.3339
$ws.Range("M59").Value = -10708.5
$ws.Range("N59").Value = -23613.3339

$ws.Range("H94").Value = 7298.5
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H138").Value = 3348.7
$ws.Range("I138").Value = 3348.7
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 10046.1
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -4906.099999999999
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4059.2307
$ws.Range("I122").Value = 2661
$ws.Range("J122").Value = 5690.5
$ws.Range("K122").Value = 7983
$ws.Range("L122").Value = 17071.5
$ws.Range("M122").Value = -5533
$ws.Range("N122").Value = -21971.5

$ws.Range("H132").Value = 3101.8
$ws.Range("I132").Value = 2836.6667
$ws.Range("K132").Value = 8510.000100000001
$ws.Range("M132").Value = -5980.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4289.125
$ws.Range("I132").Value = 4157.3
$ws.Range("J132").Value = 4948.25
$ws.Range("K132").Value = 12471.9
$ws.Range("L132").Value = 14844.75
$ws.Range("M132").Value = -9941.900000000001
$ws.Range("N132").Value = -19904.75

$ws.Range("H136").Value = 2681.125
$ws.Range("J136").Value = 1994
$ws.Range("L136").Value = 5982
$ws.Range("N136").Value = -11082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9173
$ws.Range("J45").Value = 9173
$ws.Range("L45").Value = 9173
$ws.Range("N45").Value = -10155

$ws.Range("H122").Value = 3766.4517
$ws.Range("I122").Value = 3792.476
$ws.Range("K122").Value = 11377.428
$ws.Range("M122").Value = -8927.428

$ws.Range("H125").Value = 39998.5
$ws.Range("J125").Value = 39998.5
$ws.Range("L125").Value = 39998.5
$ws.Range("N125").Value = -49838.5
